$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("meta")

# The "meta" sheet holds key/value pairs in columns A/B, one pair per row,
# terminated by a single blank "key-styled" cell. Row 10 currently is that
# blank terminator (A10, formatted like the key column, no value).
#
# Add a new "style" / "default" key-value pair in row 10 (reusing the key
# column's existing bold/orange formatting), and push the blank terminator
# row down to row 11 with the same formatting it had before.

# Give A10 the same format as the other key cells (copy from A9) before
# filling in the new pair.
$meta.Range("A9").Copy() | Out-Null
$meta.Range("A10").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$meta.Range("A10").Value = "style"
$meta.Range("B10").Value = "default"

# Recreate the blank terminator cell on the next row, keeping its format.
$meta.Range("A10").Copy() | Out-Null
$meta.Range("A11").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$excel.CutCopyMode = 0
